$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-02-21 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-22 Saturday", 2) | Out-Null
$d.Content.Find.Execute("13×57=741", $true, $false, $false, $false, $false, $true, 1, $false, "66×73=4818", 2) | Out-Null
$d.Content.Find.Execute("25×80=2000", $true, $false, $false, $false, $false, $true, 1, $false, "65×28=1820", 2) | Out-Null
$d.Content.Find.Execute("22×20=440", $true, $false, $false, $false, $false, $true, 1, $false, "66×29=1914", 2) | Out-Null
$d.Content.Find.Execute("65×86=5590", $true, $false, $false, $false, $false, $true, 1, $false, "97×76=7372", 2) | Out-Null
$d.Content.Find.Execute("97×85=8245", $true, $false, $false, $false, $false, $true, 1, $false, "89×40=3560", 2) | Out-Null
$d.Content.Find.Execute("77×92=7084", $true, $false, $false, $false, $false, $true, 1, $false, "52×80=4160", 2) | Out-Null
$d.Content.Find.Execute("32×13=416", $true, $false, $false, $false, $false, $true, 1, $false, "64×68=4352", 2) | Out-Null
$d.Content.Find.Execute("65×21=1365", $true, $false, $false, $false, $false, $true, 1, $false, "54×25=1350", 2) | Out-Null
$d.Content.Find.Execute("63×22=1386", $true, $false, $false, $false, $false, $true, 1, $false, "53×46=2438", 2) | Out-Null
$d.Content.Find.Execute("74×94=6956", $true, $false, $false, $false, $false, $true, 1, $false, "35×11=385", 2) | Out-Null
$d.Content.Find.Execute("28×23=644", $true, $false, $false, $false, $false, $true, 1, $false, "52×62=3224", 2) | Out-Null
$d.Content.Find.Execute("24×14=336", $true, $false, $false, $false, $false, $true, 1, $false, "39×60=2340", 2) | Out-Null
$d.Content.Find.Execute("29×50=1450", $true, $false, $false, $false, $false, $true, 1, $false, "37×37=1369", 2) | Out-Null
$d.Content.Find.Execute("67×63=4221", $true, $false, $false, $false, $false, $true, 1, $false, "78×91=7098", 2) | Out-Null
$d.Content.Find.Execute("79×62=4898", $true, $false, $false, $false, $false, $true, 1, $false, "76×36=2736", 2) | Out-Null
$d.Content.Find.Execute("38×17=646", $true, $false, $false, $false, $false, $true, 1, $false, "89×18=1602", 2) | Out-Null
$d.Content.Find.Execute("62×78=4836", $true, $false, $false, $false, $false, $true, 1, $false, "47×97=4559", 2) | Out-Null
$d.Content.Find.Execute("74×22=1628", $true, $false, $false, $false, $false, $true, 1, $false, "78×75=5850", 2) | Out-Null
$d.Content.Find.Execute("24×62=1488", $true, $false, $false, $false, $false, $true, 1, $false, "46×83=3818", 2) | Out-Null
$d.Content.Find.Execute("38×80=3040", $true, $false, $false, $false, $false, $true, 1, $false, "41×17=697", 2) | Out-Null
$d.Content.Find.Execute("13×44=572", $true, $false, $false, $false, $false, $true, 1, $false, "82×30=2460", 2) | Out-Null
$d.Content.Find.Execute("85×77=6545", $true, $false, $false, $false, $false, $true, 1, $false, "33×93=3069", 2) | Out-Null
$d.Content.Find.Execute("70×80=5600", $true, $false, $false, $false, $false, $true, 1, $false, "65×32=2080", 2) | Out-Null
$d.Content.Find.Execute("72×78=5616", $true, $false, $false, $false, $false, $true, 1, $false, "36×95=3420", 2) | Out-Null
$d.Content.Find.Execute("74×65=4810", $true, $false, $false, $false, $false, $true, 1, $false, "89×73=6497", 2) | Out-Null
